$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.890.48'
$ws.Range('E2').Value = '  -2.20%  '
$ws.Range('D3').Value = '1.865.07'
$ws.Range('E3').Value = '  -2.52%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '''311.73'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = '''0.4961'
$ws.Range('E7').Value = '  -3.66%  '
$ws.Range('D8').Value = '''0.3793'
$ws.Range('E8').Value = '  -4.58%  '
$ws.Range('D9').Value = '''0.08881'
$ws.Range('E9').Value = '  -9.42%  '
$ws.Range('D10').Value = '''1.113'
$ws.Range('E10').Value = '  -3.30%  '
$ws.Range('D11').Value = '''41.45'
$ws.Range('E11').Value = '  -1.85%  '
$ws.Range('D12').Value = '''6.291'
$ws.Range('E12').Value = '  -3.68%  '
$ws.Range('D13').Value = '''20.56'
$ws.Range('E13').Value = '  -2.91%  '
$ws.Range('D14').Value = '1.859.53'
$ws.Range('E14').Value = '  -2.68%  '
$ws.Range('E15').Value = '  -3.86%  '
$ws.Range('E17').Value = '  -3.75%  '
$ws.Range('E18').Value = '  -4.39%  '
$ws.Range('D19').Value = '''0.06626'
$ws.Range('D20').Value = '''17.81'
$ws.Range('E20').Value = '  -2.59%  '
$ws.Range('E21').Value = '  +0.15%  '
$ws.Range('D22').Value = '''6.077'
$ws.Range('E22').Value = '  -3.89%  '
$ws.Range('D23').Value = '27.927.37'
$ws.Range('E23').Value = '  -2.28%  '
$ws.Range('D24').Value = '''11.34'
$ws.Range('E24').Value = '  -1.61%  '
$ws.Range('D25').Value = '''2.283'
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('D26').Value = '2.084.62'
$ws.Range('E26').Value = '  -2.10%  '
$ws.Range('D27').Value = '''2.506'
$ws.Range('E27').Value = '  -6.54%  '
$ws.Range('D28').Value = '''157.93'
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('E29').Value = '  -3.14%  '
$ws.Range('D30').Value = '''125.67'
$ws.Range('E30').Value = '  -2.69%  '
$ws.Range('D31').Value = '''0.1052'
$ws.Range('E31').Value = '  -2.37%  '
$ws.Range('E32').Value = '  -5.74%  '
$ws.Range('D33').Value = '''5.559'
$ws.Range('E33').Value = '  -3.26%  '
$ws.Range('D34').Value = '''3.585'
$ws.Range('E34').Value = '  -1.26%  '
$ws.Range('D35').Value = '''9.288'
$ws.Range('E35').Value = '  -5.84%  '
$ws.Range('E36').Value = '  -4.21%  '
$ws.Range('D37').Value = '''0.02393'
$ws.Range('E37').Value = '  -1.90%  '
$ws.Range('D38').Value = '''0.2172'
$ws.Range('E38').Value = '  -2.06%  '
$ws.Range('D39').Value = '''1.265'
$ws.Range('E39').Value = '  +6.56%  '
$ws.Range('D40').Value = '''1.195'
$ws.Range('E40').Value = '  -6.42%  '
$ws.Range('D41').Value = '''11.61'
$ws.Range('E41').Value = '  -1.63%  '
$ws.Range('D42').Value = '''0.6330'
$ws.Range('E42').Value = '  -2.40%  '
$ws.Range('D43').Value = '''4.874'
$ws.Range('E43').Value = '  -4.33%  '
$ws.Range('D44').Value = '''1.001'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''13.15'
$ws.Range('E45').Value = '  -3.02%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '''0.5955'
$ws.Range('E46').Value = '  -2.50%  '
$ws.Range('D47').Value = '''1.279'
$ws.Range('E47').Value = '  -0.85%  '
$ws.Range('D48').Value = '''3.666'
$ws.Range('E48').Value = '  -2.76%  '
$ws.Range('D49').Value = '''1.208'
$ws.Range('E49').Value = '  +0.30%  '
$ws.Range('D50').Value = '''1.958'
$ws.Range('E50').Value = '  -4.35%  '
$ws.Range('D51').Value = '''120.53'
$ws.Range('E51').Value = '  -3.42%  '

Write-Output "Applied 93 cell updates"
